$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Harvesting_Isolation")
Write-Host "ok"
